$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.182.31'
$ws.Range("E2").Value = '  -2.07%  '

# Row 3
$ws.Range("D3").Value = '3.680.59'
$ws.Range("E3").Value = '  -3.05%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '681.82'
$ws.Range("E5").Value = '  -3.62%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.55'
$ws.Range("E6").Value = '  -4.70%  '

# Row 7
$ws.Range("D7").Value = '3.679.57'
$ws.Range("E7").Value = '  -3.04%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.498'
$ws.Range("E9").Value = '  -4.49%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.149'
$ws.Range("E10").Value = '  -7.60%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.21'
$ws.Range("E11").Value = '  -2.11%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.451'
$ws.Range("E12").Value = '  -1.23%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000236'
$ws.Range("E13").Value = '  -7.08%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.50'
$ws.Range("E14").Value = '  -7.39%  '

# Row 15
$ws.Range("D15").Value = '4.304.83'
$ws.Range("E15").Value = '  -2.89%  '

# Row 16
$ws.Range("D16").Value = '3.681.70'
$ws.Range("E16").Value = '  -3.53%  '

# Row 17
$ws.Range("D17").Value = '69.272.14'
$ws.Range("E17").Value = '  -1.97%  '

# Row 18
$ws.Range("E18").Value = '  -1.78%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.33'
$ws.Range("E19").Value = '  -6.01%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.65'
$ws.Range("E20").Value = '  -6.90%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '483.28'
$ws.Range("E21").Value = '  -2.18%  '

# Row 22
$ws.Range("E22").Value = '  -7.66%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.665'
$ws.Range("E23").Value = '  -8.65%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.74'
$ws.Range("E24").Value = '  -6.12%  '

# Row 25
$ws.Range("D25").Value = '3.826.14'
$ws.Range("E25").Value = '  -2.99%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.55'
$ws.Range("E26").Value = '  -4.38%  '

# Row 27
$ws.Range("E27").Value = '  -12.11%  '

# Row 28
$ws.Range("E28").Value = '  -0.07%  '

# Row 29
$ws.Range("E29").Value = '  -8.82%  '

# Row 30
$ws.Range("E30").Value = '  -9.81%  '

# Row 31
$ws.Range("E31").Value = '  -11.35%  '

# Row 32
$ws.Range("E32").Value = '  -4.94%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.76'
$ws.Range("E33").Value = '  -7.87%  '

# Row 34
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.98'
$ws.Range("E34").Value = '  -7.34%  '

# Row 35
$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.02%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.163'
$ws.Range("E36").Value = '  -6.51%  '

# Row 37
$ws.Range("D37").Value = '3.649.04'
$ws.Range("E37").Value = '  -3.05%  '

# Row 38
$ws.Range("E38").Value = '  -5.97%  '

# Row 39
$ws.Range("E39").Value = '  +2.36%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0945'
$ws.Range("E40").Value = '  -6.73%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.17'
$ws.Range("E42").Value = '  -5.82%  '

# Row 43
$ws.Range("E43").Value = '  +0.01%  '

# Row 44
$ws.Range("E44").Value = '  -7.83%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '157.00'
$ws.Range("E45").Value = '  -4.51%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '48.00'
$ws.Range("E46").Value = '  -1.81%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.78'
$ws.Range("E47").Value = '  -15.58%  '

# Row 48
$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '392.99'
$ws.Range("E48").Value = '  -6.90%  '

# Row 49
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.31'
$ws.Range("E49").Value = '  -4.06%  '

# Row 50
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000278'
$ws.Range("E50").Value = '  -13.16%  '

# Row 51
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.08'
$ws.Range("E51").Value = '  -7.11%  '
